$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 48
$ws.Range("F2").Value = 28
$ws.Range("H2").Value = 28

# Row 5
$ws.Range("E5").Value = 4

# Row 17
$ws.Range("E17").Value = 113
$ws.Range("F17").Value = 54
$ws.Range("H17").Value = 54

# Row 18
$ws.Range("E18").Value = 106

# Row 33
$ws.Range("E33").Value = 40

# Row 34
$ws.Range("E34").Value = 20

# Row 42
$ws.Range("E42").Value = 35

# Row 45
$ws.Range("E45").Value = 26
$ws.Range("F45").Value = 14
$ws.Range("H45").Value = 14

# Row 49
$ws.Range("E49").Value = 67
$ws.Range("F49").Value = 39
$ws.Range("H49").Value = 39

# Row 57
$ws.Range("E57").Value = 14

# Row 64
$ws.Range("E64").Value = 34

# Row 68
$ws.Range("E68").Value = 16

# Row 71
$ws.Range("E71").Value = 34
$ws.Range("F71").Value = 16
$ws.Range("H71").Value = 16

# Row 79
$ws.Range("E79").Value = 38
$ws.Range("F79").Value = 18
$ws.Range("H79").Value = 18

# Row 89
$ws.Range("E89").Value = 34
